$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price column (D) - force text format to preserve exact string representation
$priceCells = @{
    "D2" = "42.722.34"
    "D3" = "2.533.93"
    "D5" = "317.34"
    "D6" = "95.44"
    "D7" = "0.578"
    "D9" = "0.529"
    "D10" = "36.18"
    "D11" = "0.0809"
    "D12" = "7.57"
    "D14" = "2.921.79"
    "D15" = "15.44"
    "D16" = "2.524.82"
    "D17" = "0.850"
    "D18" = "42.693.58"
    "D19" = "13.03"
    "D20" = "6.58"
    "D21" = "0.0₃0962"
    "D22" = "70.20"
    "D23" = "251.37"
    "D24" = "2.98"
    "D26" = "26.71"
    "D28" = "2.41"
    "D29" = "39.15"
    "D30" = "10.16"
    "D31" = "6.08"
    "D32" = "154.90"
    "D33" = "2.13"
    "D34" = "19.01"
    "D36" = "0.0789"
    "D39" = "0.119"
    "D40" = "23.53"
    "D42" = "3.81"
    "D44" = "0.0301"
    "D46" = "2.018.55"
    "D47" = "85.61"
    "D48" = "8.82"
    "D49" = "2.776.21"
    "D50" = "74.14"
    "D51" = "102.79"
}
foreach ($addr in $priceCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceCells[$addr]
    $cell.Style = "Normal"
}

# Update Volume(1h) column (E)
$volumeCells = @{
    "E2" = "  -0.63%  "
    "E3" = "  -0.09%  "
    "E4" = "  +0.00%  "
    "E5" = "  +3.50%  "
    "E6" = "  -5.38%  "
    "E7" = "  -1.12%  "
    "E8" = "  -0.07%  "
    "E9" = "  -3.35%  "
    "E10" = "  -2.90%  "
    "E11" = "  -1.02%  "
    "E12" = "  -1.04%  "
    "E13" = "  -0.18%  "
    "E14" = "  +0.05%  "
    "E15" = "  +1.35%  "
    "E16" = "  -1.87%  "
    "E17" = "  -2.19%  "
    "E18" = "  -0.63%  "
    "E19" = "  -0.62%  "
    "E20" = "  +1.12%  "
    "E22" = "  -2.14%  "
    "E23" = "  -1.24%  "
    "E24" = "  +1.18%  "
    "E25" = "  -2.58%  "
    "E26" = "  -2.45%  "
    "E27" = "  +0.05%  "
    "E28" = "  +3.35%  "
    "E29" = "  +0.28%  "
    "E30" = "  -3.74%  "
    "E31" = "  -1.70%  "
    "E32" = "  -2.44%  "
    "E33" = "  +0.58%  "
    "E34" = "  +2.62%  "
    "E35" = "  -1.42%  "
    "E36" = "  -1.15%  "
    "E37" = "  -0.45%  "
    "E38" = "  -3.63%  "
    "E39" = "  -1.20%  "
    "E40" = "  -3.50%  "
    "E42" = "  -2.85%  "
    "E43" = "  +0.29%  "
    "E44" = "  -1.24%  "
    "E45" = "  -5.16%  "
    "E46" = "  -1.31%  "
    "E47" = "  -0.70%  "
    "E48" = "  -2.04%  "
    "E49" = "  -0.20%  "
    "E50" = "  +0.86%  "
    "E51" = "  -0.87%  "
}
foreach ($addr in $volumeCells.Keys) {
    $ws.Range($addr).Value = $volumeCells[$addr]
}
